# Update automation script and BrowserStack automation script video url
# - Sheet1 (xl/worksheets/sheet2.xml): swap RunStatus for rows 30/31 (D30<->D31)
# - Insert two new steps (rows 45 & 46) before the final "quit" step (old row 45,
#   which becomes row 47), adding two new xpath locators to the shared strings
#   table, and update the selected cell to reflect the new last-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Swap the RunStatus (column D) values on rows 30 and 31 -------------
$ws.Range("D30").Value = "N"
$ws.Range("D31").Value = "Y"

# --- 2. Insert two fresh rows right before the old row 45 ("quit" step) ----
# This pushes the old row 45 down to row 47 and leaves two blank rows
# (45 and 46) ready to be filled in with the new steps.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# --- 3. Populate new row 45: click the user-dropdown -----------------------
$ws.Range("A45").Value = "'38"
$ws.Range("B45").Value = "IshinePortal"
$ws.Range("C45").Value = "IshineOTPField"
$ws.Range("D45").Value = "Y"
$ws.Range("E45").Value = "xpath"
$ws.Range("F45").Value = "(//a[@id='user-dropdown'])[1]"
$ws.Range("H45").Value = "click"
$ws.Range("J45").Value = "TC_01_05"
$ws.Range("K45").Value = "User should be able to login after entering OTP"
$ws.Range("L45").Value = "SC_38"

# --- 4. Populate new row 46: click the dropdown's logout item --------------
$ws.Range("A46").Value = "'38"
$ws.Range("B46").Value = "IshinePortal"
$ws.Range("C46").Value = "IshineOTPField"
$ws.Range("D46").Value = "Y"
$ws.Range("E46").Value = "xpath"
$ws.Range("F46").Value = "(//a[@class = 'dropdown-item'])[2]"
$ws.Range("H46").Value = "click"
$ws.Range("J46").Value = "TC_01_05"
$ws.Range("K46").Value = "User should be able to login after entering OTP"
$ws.Range("L46").Value = "SC_38"

# --- 5. Row 47 (old row 45, "quit" step) keeps its original content; it was
# shifted down automatically by the Insert() calls above, so nothing else to
# change there.

# --- 6. Reflect the new selection / scroll position in the sheet view ------
[void]$ws.Range("F46").Select()
